$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-13 is updated to the new serial date value (2023-09-14).
$ws.Range("C2:C13").Value = 45183
